$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact literal text value (preserves leading zeros,
# trailing zeros, percent signs, etc. exactly as authored) without Excel
# re-interpreting / reformatting the string as a number.
function Set-TextValue($ws, [string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '285.26'
Set-TextValue $ws 'E2' '2.24%'

# Row 3
Set-TextValue $ws 'D3' '28.70'
Set-TextValue $ws 'E3' '4.79%'

# Row 4
Set-TextValue $ws 'D4' '4.986'
Set-TextValue $ws 'E4' '3.09%'

# Row 5
Set-TextValue $ws 'D5' '0.06701'
Set-TextValue $ws 'E5' '5.11%'

# Row 6
Set-TextValue $ws 'D6' '7.347'
Set-TextValue $ws 'E6' '4.46%'

# Row 7
Set-TextValue $ws 'B7' 'FTXToken'
Set-TextValue $ws 'C7' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws 'D7' '1.359'
Set-TextValue $ws 'E7' '5.81%'

# Row 8
Set-TextValue $ws 'B8' 'MXToken'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D8' '0.9395'
Set-TextValue $ws 'E8' '5.28%'

# Row 9
Set-TextValue $ws 'B9' 'WazirX'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D9' '0.1575'
Set-TextValue $ws 'E9' '3.30%'

# Row 10
Set-TextValue $ws 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D10' '0.06639'
Set-TextValue $ws 'E10' '14.77%'

# Row 11
Set-TextValue $ws 'B11' 'MandalaExchangeToken'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.07639'
Set-TextValue $ws 'E11' '1.92%'

# Row 12
Set-TextValue $ws 'B12' 'BitrueCoin'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D12' '0.02922'
Set-TextValue $ws 'E12' '0.15%'

# Row 13
Set-TextValue $ws 'B13' 'BitMartToken'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D13' '0.08991'
Set-TextValue $ws 'E13' '0.04%'

# Row 14
Set-TextValue $ws 'B14' 'BitForexToken'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D14' '0.001581'
Set-TextValue $ws 'E14' '0.62%'

# Row 15
Set-TextValue $ws 'B15' 'CoinExToken'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws 'D15' '0.04473'
Set-TextValue $ws 'E15' '1.70%'

# Row 16
Set-TextValue $ws 'B16' 'One'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws 'D16' '0.0006433'
Set-TextValue $ws 'E16' '0.54%'

# Row 17
Set-TextValue $ws 'B17' 'TigerCash'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D17' '0.006535'
Set-TextValue $ws 'E17' '7.01%'

# Row 18
Set-TextValue $ws 'B18' 'LEO'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D18' '3.489'
Set-TextValue $ws 'E18' '0.47%'

# Row 19
Set-TextValue $ws 'B19' 'GateToken'
Set-TextValue $ws 'C19' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D19' '3.379'
Set-TextValue $ws 'E19' '2.29%'

# Row 20
Set-TextValue $ws 'D20' '2.235'
Set-TextValue $ws 'E20' '-2.70%'

# Row 22
Set-TextValue $ws 'E22' '-3.13%'

# Row 23
Set-TextValue $ws 'D23' '4.048'

# Row 24
Set-TextValue $ws 'D24' '0.1521'
Set-TextValue $ws 'E24' '1.17%'

# Row 25
Set-TextValue $ws 'D25' '0.001177'
Set-TextValue $ws 'E25' '0.13%'

# Row 26
Set-TextValue $ws 'D26' '0.004487'
Set-TextValue $ws 'E26' '4.80%'

# Row 27
Set-TextValue $ws 'D27' '0.0001245'
Set-TextValue $ws 'E27' '5.60%'

# Row 28
Set-TextValue $ws 'D28' '0.0001612'
Set-TextValue $ws 'E28' '-2.35%'

# Row 40
Set-TextValue $ws 'D40' '0.04194'
Set-TextValue $ws 'E40' '3.67%'

# Row 41
Set-TextValue $ws 'D41' '0.006729'
Set-TextValue $ws 'E41' '1.42%'

# Row 42
Set-TextValue $ws 'D42' '0.1251'
Set-TextValue $ws 'E42' '-10.93%'

# Row 43
Set-TextValue $ws 'D43' '0.002013'
Set-TextValue $ws 'E43' '-2.25%'

# Row 44
Set-TextValue $ws 'D44' '0.01218'
Set-TextValue $ws 'E44' '9.30%'

# Row 45
Set-TextValue $ws 'D45' '0.00005594'
Set-TextValue $ws 'E45' '1.03%'

# Row 46
Set-TextValue $ws 'E46' '25.93%'

# Row 47
Set-TextValue $ws 'D47' '0.01302'
Set-TextValue $ws 'E47' '-29.48%'
